$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A1").Value = "  No aplicable"
$ws.Range("A2").Value = "1951 a 1960"
$ws.Range("A3").Value = "1941 - 1950"
$ws.Range("A4").Value = "Antes 1900"
$ws.Range("A6").Value = "1900 - 1920"
$ws.Range("A7").Value = "1971 a 1980"
$ws.Range("A9").Value = "1961 - 1970"
$ws.Range("A13").Value = "1921 a 1940"
$ws.Range("A22").Value = "1981 a 1990"
$ws.Range("A28").Value = "1951 - 1960"
$ws.Range("A29").Value = "1991 a 2001"
$ws.Range("A31").Value = "1941 a 1950"
$ws.Range("A32").Value = "1921 - 1940"
$ws.Range("A35").Value = "1971 - 1980"
$ws.Range("A36").Value = "1900 a 1920"
$ws.Range("A37").Value = "1981 - 1990"
$ws.Range("A38").Value = "1961 a 1970"
$ws.Range("A39").Value = " Antes de 1900"
